# Insert a new data row at row 620 (a new daily price observation for
# Coliflor at "Vega Modelo de Temuco"), shifting all the existing rows
# 620:693 down to 621:694 -- matches the dimension change A1:R693 -> A1:R694.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 620..693 down to 621..694, opening up a blank row at 620.
$ws.Rows("620:620").Insert()

# Populate the newly-opened row 620 with the new observation.
$ws.Range("A620").Value = 10
$ws.Range("B620").Value = "Vega Modelo de Temuco"
$ws.Range("C620").Value = "La Araucanía"
$ws.Range("D620").Value = 45142
$ws.Range("E620").Value = 9
$ws.Range("F620").Value = 100112008
$ws.Range("G620").Value = "Coliflor"
$ws.Range("H620").Value = "Sin especificar"
$ws.Range("I620").Value = "Primera"
$ws.Range("J620").Value = 600
$ws.Range("K620").Value = 1000
$ws.Range("L620").Value = 1000
$ws.Range("M620").Value = 1000
$ws.Range("N620").Value = "$/unidad"
$ws.Range("O620").Value = "Región Metropolitana"
$ws.Range("P620").Value = 1000
$ws.Range("Q620").Value = 1
$ws.Range("R620").Value = "Hortaliza"
